# Ridge Classifier ML documentation update:
# - Splits the old single "Drop" column into two explicit boolean columns:
#     "Drop Medical Specialty" (col E) and "Drop Emergency" (col F), and
#   moves "Score" to col G.
# - Back-fills Drop Medical Specialty = Y for existing rows 2-7 (their
#   historical runs all dropped Medical Specialty) and Drop Emergency = N.
# - Corrects the row 6 score.
# - Appends two new experiment rows (8-9) for the Drop Emergency variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Grab format templates from pristine cells before anything else in
#    the sheet gets rearranged:
#      - E4 currently carries the "Var" font / left-center alignment
#        style (s=1) that the new Score column (G) should keep using.
#      - A7 currently carries the bold "last row" style (s=2) that the
#        two newly appended rows should use for their text columns.
# ---------------------------------------------------------------------
$ws.Range("E4").Copy()
$ws.Range("G4:G9").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("A8:F8").PasteSpecial(-4122)
$ws.Range("A9:E9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 1. Header row: insert the two new columns, shift Score to G, drop H1.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "Drop Medical Specialty"
$ws.Range("F1").Value = "Drop Emergency"
$ws.Range("G1").Value = "Score"
$ws.Range("H1").Value = ""

# ---------------------------------------------------------------------
# 2. Existing data rows 2-6: shift the Score value from E to G and fill
#    in the two new boolean columns (all dropped Medical Specialty,
#    none dropped Emergency). Rows 4-6's old E cell had the special
#    "Var" font style that now belongs on G, so it's reset to plain
#    formatting once its value/style have been copied out.
# ---------------------------------------------------------------------
$ws.Range("G2").Value = $ws.Range("E2").Value()
$ws.Range("E2").Value = "Y"
$ws.Range("F2").Value = "N"

$ws.Range("G3").Value = $ws.Range("E3").Value()
$ws.Range("E3").Value = "Y"
$ws.Range("F3").Value = "N"

$ws.Range("G4").Value = 0.63963602301619105
$ws.Range("E4").ClearFormats()
$ws.Range("E4").Value = "Y"
$ws.Range("F4").Value = "N"

$ws.Range("G5").Value = 0.64993978321959001
$ws.Range("E5").ClearFormats()
$ws.Range("E5").Value = "Y"
$ws.Range("F5").Value = "N"

# Row 6 also got a corrected score while the columns were reshuffled.
$ws.Range("G6").Value = 0.64077345109059203
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = "Y"
$ws.Range("F6").Value = "N"

# ---------------------------------------------------------------------
# 3. Row 7 (previously the last/bold row): shift score to G, fill new
#    columns (E7/F7 already carry the bold "last row" style copied from
#    A7 above), and add the empty bold H7 placeholder that keeps the
#    dimension/formatting consistent with the new last bold row.
# ---------------------------------------------------------------------
$ws.Range("G7").Value = 0.65067576609126099
$ws.Range("E7").Value = "Y"
$ws.Range("F7").Value = "N"
$ws.Range("H7").Font.Bold = $true

# ---------------------------------------------------------------------
# 4. New row 8: Ridge Classifier, Y, Y, N, N (Drop Medical Specialty),
#    N (Drop Emergency), score 0.64893617021276595. A8:F8 already carry
#    the bold "last row" style copied from A7 above.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Ridge Classifier"
$ws.Range("B8").Value = "Y"
$ws.Range("C8").Value = "Y"
$ws.Range("D8").Value = "N"
$ws.Range("E8").Value = "N"
$ws.Range("F8").Value = "N"
$ws.Range("G8").Value = 0.64893617021276595
$ws.Range("H8").Font.Bold = $true

# ---------------------------------------------------------------------
# 5. New row 9: Ridge Classifier, Y, Y, Y, N (Drop Medical Specialty),
#    Y (Drop Emergency), score 0.64732650739476605 -- note row 9 has no
#    H cell. F9 keeps the "Var" font style (like the Score column)
#    rather than the plain bold row style used by the rest of row 9.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Ridge Classifier"
$ws.Range("B9").Value = "Y"
$ws.Range("C9").Value = "Y"
$ws.Range("D9").Value = "Y"
$ws.Range("E9").Value = "N"
$ws.Range("G9").Value = 0.64732650739476605

$ws.Range("E4").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = "Y"

# ---------------------------------------------------------------------
# 6. Column widths for the two new columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 13.57
$ws.Columns.Item(6).ColumnWidth = 13.57

# ---------------------------------------------------------------------
# 7. Selection / active cell moves below the new last row.
# ---------------------------------------------------------------------
$ws.Range("A10").Select()

Write-Host "done"
